$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (this shifts all existing data rows down by one,
# row 2 -> row 3, row 3 -> row 4, ..., row 40 -> row 41), matching the
# diff where every existing record shifted down and a brand-new record
# (previously absent) was introduced at the top as the new row 2.
$ws.Rows.Item(2).Insert()

# The inserted row copies formatting (bold header style) from the row
# above it; clear that so the new row matches the plain data-row styling
# used throughout the rest of the sheet.
$ws.Range("A2:R2").ClearFormats()

# The Fecha (date) column uses a custom date/time number format on every
# data row; re-apply it to the new row's D cell before writing the value.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with the new record's data.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C2").Value = 'Arica y Parinacota'
$ws.Range("D2").Value = 44860
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 100112028
$ws.Range("G2").Value = 'Sandia'
$ws.Range("H2").Value = 'Sin especificar'
$ws.Range("I2").Value = 'Segunda'
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 730
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = 739
$ws.Range("N2").Value = '$/kilo (volumen en unidades)'
$ws.Range("O2").Value = 'Perú'
$ws.Range("P2").Value = 739
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 'Hortaliza'
